$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "50.995.51"
$ws.Range("E2").Value = "  -1.22%  "

$ws.Range("D3").Value = "2.944.76"
$ws.Range("E3").Value = "  -1.40%  "

$ws.Range("E4").Value = "  -0.03%  "

$ws.Range("D5").Value = "'375.26"
$ws.Range("E5").Value = "  -2.44%  "

$ws.Range("D6").Value = "'101.13"
$ws.Range("E6").Value = "  -3.36%  "

$ws.Range("E7").Value = "  -1.31%  "

$ws.Range("E8").Value = "  -0.02%  "

$ws.Range("D9").Value = "'0.590"
$ws.Range("E9").Value = "  -0.99%  "

$ws.Range("D10").Value = "'36.40"
$ws.Range("E10").Value = "  -2.31%  "

$ws.Range("E11").Value = "  -0.73%  "

$ws.Range("D12").Value = "'0.0851"
$ws.Range("E12").Value = "  +0.36%  "

$ws.Range("D13").Value = "3.413.18"
$ws.Range("E13").Value = "  -1.52%  "

$ws.Range("E14").Value = "  -1.23%  "

$ws.Range("D15").Value = "'7.59"
$ws.Range("E15").Value = "  -0.34%  "

$ws.Range("D16").Value = "'11.24"
$ws.Range("E16").Value = "  +50.59%  "

$ws.Range("D17").Value = "2.950.28"
$ws.Range("E17").Value = "  -1.42%  "

$ws.Range("D18").Value = "'0.999"
$ws.Range("E18").Value = "  -0.71%  "

$ws.Range("D19").Value = "50.969.54"
$ws.Range("E19").Value = "  -1.28%  "

$ws.Range("E20").Value = "  -5.90%  "

$ws.Range("D21").Value = "'12.48"
$ws.Range("E21").Value = "  -3.40%  "

$ws.Range("D22").Value = "0.0₃0955"
$ws.Range("E22").Value = "  -1.04%  "

$ws.Range("D23").Value = "'265.68"
$ws.Range("E23").Value = "  +0.78%  "

$ws.Range("D24").Value = "'68.82"
$ws.Range("E24").Value = "  -0.51%  "

$ws.Range("D25").Value = "'3.14"
$ws.Range("E25").Value = "  +7.14%  "

$ws.Range("D26").Value = "'8.13"
$ws.Range("E26").Value = "  -2.73%  "

$ws.Range("D27").Value = "'7.57"
$ws.Range("E27").Value = "  -1.89%  "

$ws.Range("E28").Value = "  -0.25%  "

$ws.Range("E29").Value = "  +0.02%  "

$ws.Range("B30").Value = "Kaspa"
$ws.Range("C30").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D30").Value = "'0.164"
$ws.Range("E30").Value = "  -3.44%  "

$ws.Range("B31").Value = "EthereumClassic"
$ws.Range("C31").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D31").Value = "'25.70"
$ws.Range("E31").Value = "  -1.34%  "

$ws.Range("E32").Value = "  -4.79%  "

$ws.Range("D33").Value = "'10.00"
$ws.Range("E33").Value = "  +1.03%  "

$ws.Range("D34").Value = "'50.93"
$ws.Range("E34").Value = "  -0.30%  "

$ws.Range("E35").Value = "  -0.53%  "

$ws.Range("D36").Value = "'33.46"
$ws.Range("E36").Value = "  -4.21%  "

$ws.Range("D37").Value = "'0.0444"
$ws.Range("E37").Value = "  -2.09%  "

$ws.Range("E38").Value = "  -0.12%  "

$ws.Range("D39").Value = "'3.17"
$ws.Range("E39").Value = "  +4.48%  "

$ws.Range("E40").Value = "  -0.61%  "

$ws.Range("D41").Value = "'16.52"
$ws.Range("E41").Value = "  -3.26%  "

$ws.Range("E42").Value = "  -2.01%  "

$ws.Range("D43").Value = "'2.47"
$ws.Range("E43").Value = "  -5.33%  "

$ws.Range("D44").Value = "'120.66"
$ws.Range("E44").Value = "  -1.51%  "

$ws.Range("D45").Value = "'21.25"
$ws.Range("E45").Value = "  -2.42%  "

$ws.Range("E46").Value = "  +3.47%  "

$ws.Range("E47").Value = "  -0.69%  "

$ws.Range("D48").Value = "'0.273"
$ws.Range("E48").Value = "  -2.21%  "

$ws.Range("D49").Value = "'2.34"
$ws.Range("E49").Value = "  -1.56%  "

$ws.Range("D50").Value = "1.993.13"
$ws.Range("E50").Value = "  -2.28%  "

$ws.Range("D51").Value = "'0.0328"
$ws.Range("E51").Value = "  -1.41%  "
